$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '60.978.52'
$ws.Range("E2").Value2 = '  +1.00%  '
$ws.Range("D3").Value2 = '3.377.12'
$ws.Range("E3").Value2 = '  +0.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  +0.12%  '
$ws.Range("E7").Value2 = '  +0.02%  '
$ws.Range("E8").Value2 = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.62'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = '  +1.57%  '
$ws.Range("E10").Value2 = '  -1.22%  '
$ws.Range("E11").Value2 = '  -0.21%  '
$ws.Range("D12").Value2 = '3.954.43'
$ws.Range("E12").Value2 = '  +0.19%  '
$ws.Range("E13").Value2 = '  +2.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = '  -0.89%  '
$ws.Range("D15").Value2 = '3.373.68'
$ws.Range("E15").Value2 = '  -0.04%  '
$ws.Range("E16").Value2 = '  -0.27%  '
$ws.Range("D17").Value2 = '61.078.00'
$ws.Range("E17").Value2 = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = '  -1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  -2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = '  -1.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  -1.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  +3.19%  '
$ws.Range("E23").Value2 = '  -1.03%  '
$ws.Range("E24").Value2 = '  +0.15%  '
$ws.Range("B25").Value2 = 'PEPE'
$ws.Range("C25").Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000113'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = '  -1.30%  '
$ws.Range("B26").Value2 = 'WrappedeETH'
$ws.Range("C26").Value2 = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value2 = '3.515.30'
$ws.Range("E26").Value2 = '  -0.08%  '
$ws.Range("E27").Value2 = '  +7.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = '  -2.42%  '
$ws.Range("E30").Value2 = '  +0.36%  '
$ws.Range("E31").Value2 = '  -0.31%  '
$ws.Range("E33").Value2 = '  -3.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = '  -1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = '  +0.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '166.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = '  -1.18%  '
$ws.Range("D37").Value2 = '3.412.99'
$ws.Range("E37").Value2 = '  +0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = '  +0.91%  '
$ws.Range("E39").Value2 = '  -2.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0762'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  -0.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = '  -4.33%  '
$ws.Range("E42").Value2 = '  -0.02%  '
$ws.Range("E43").Value2 = '  +0.30%  '
$ws.Range("E44").Value2 = '  -1.62%  '
$ws.Range("E45").Value2 = '  -3.26%  '
$ws.Range("E46").Value2 = '  -0.36%  '
$ws.Range("D47").Value2 = '2.437.82'
$ws.Range("E47").Value2 = '  -2.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = '  -1.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = '  -1.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0259'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = '  -2.95%  '
$ws.Range("E51").Value2 = '  +6.03%  '
